# Rename the three worksheets from Russian default names to English ones,
# and fix the date in the report header on the first sheet.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "Worksheet 1"
$wb.Worksheets.Item(2).Name = "Worksheet 2"
$wb.Worksheets.Item(3).Name = "Worksheet 3"

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Date: 05-10-2018 - Department: Sales department"
